# Updated cryptos list values (prices & 1h volume changes) per upstream diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking text (e.g. "212.78"); force text format so
# Excel COM does not silently coerce these assignments into real numbers,
# matching the original inline-string (text) cell type.
$ws.Range('D2:D51').NumberFormat = '@'

$ws.Range('D2').Value = '29.642.59'
$ws.Range('E2').Value = '  +3.64%  '

$ws.Range('D3').Value = '1.608.89'
$ws.Range('E3').Value = '  +2.88%  '

$ws.Range('E4').Value = '  +0.02%  '

$ws.Range('D5').Value = '212.78'
$ws.Range('E5').Value = '  +1.19%  '

$ws.Range('D6').Value = '0.521'
$ws.Range('E6').Value = '  +3.14%  '

$ws.Range('E7').Value = '  +0.05%  '

$ws.Range('D8').Value = '26.99'
$ws.Range('E8').Value = '  +8.50%  '

$ws.Range('D9').Value = '43.60'
$ws.Range('E9').Value = '  -1.18%  '

$ws.Range('E10').Value = '  +2.71%  '

$ws.Range('D11').Value = '0.0602'
$ws.Range('E11').Value = '  +2.64%  '

$ws.Range('D12').Value = '0.0910'
$ws.Range('E12').Value = '  +1.56%  '

$ws.Range('D13').Value = '1.837.99'
$ws.Range('E13').Value = '  +2.82%  '

$ws.Range('D14').Value = '1.607.89'
$ws.Range('E14').Value = '  +2.73%  '

$ws.Range('D15').Value = '29.654.32'
$ws.Range('E15').Value = '  +3.54%  '

$ws.Range('E16').Value = '  +4.10%  '

$ws.Range('E17').Value = '  +2.62%  '

$ws.Range('D18').Value = '63.53'
$ws.Range('E18').Value = '  +3.54%  '

$ws.Range('D19').Value = '241.12'
$ws.Range('E19').Value = '  +6.06%  '

$ws.Range('D20').Value = '7.61'
$ws.Range('E20').Value = '  +4.00%  '

$ws.Range('D21').Value = '0.0₃0695'
$ws.Range('E21').Value = '  +2.18%  '

$ws.Range('D23').Value = '4.00'
$ws.Range('E23').Value = '  +1.83%  '

$ws.Range('D24').Value = '9.25'
$ws.Range('E24').Value = '  +2.31%  '

$ws.Range('D25').Value = '2.09'
$ws.Range('E25').Value = '  +0.91%  '

$ws.Range('D26').Value = '154.88'
$ws.Range('E26').Value = '  +2.00%  '

$ws.Range('D27').Value = '15.32'
$ws.Range('E27').Value = '  +3.75%  '

$ws.Range('E28').Value = '  +3.29%  '

$ws.Range('D29').Value = '6.43'
$ws.Range('E29').Value = '  +3.06%  '

$ws.Range('E30').Value = '  +0.02%  '

$ws.Range('D31').Value = '0.0475'
$ws.Range('E31').Value = '  +3.64%  '

$ws.Range('E32').Value = '  +1.18%  '

$ws.Range('E33').Value = '  +1.58%  '

$ws.Range('E34').Value = '  +4.55%  '

$ws.Range('D35').Value = '1.420.10'
$ws.Range('E35').Value = '  +1.14%  '

$ws.Range('E36').Value = '  -0.15%  '

$ws.Range('E37').Value = '  +5.21%  '

$ws.Range('E38').Value = '  +5.60%  '

$ws.Range('E39').Value = '  +0.28%  '

$ws.Range('D40').Value = '0.0166'
$ws.Range('E40').Value = '  +2.33%  '

$ws.Range('D41').Value = '0.541'
$ws.Range('E41').Value = '  +4.60%  '

$ws.Range('E42').Value = '  +2.52%  '

$ws.Range('B43').Value = 'BitcoinSV'
$ws.Range('C43').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D43').Value = '55.10'
$ws.Range('E43').Value = '  +30.20%  '

$ws.Range('B44').Value = 'Kaspa'
$ws.Range('C44').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D44').Value = '0.0493'
$ws.Range('E44').Value = '  +6.79%  '

$ws.Range('D45').Value = '0.801'
$ws.Range('E45').Value = '  +4.49%  '

$ws.Range('D46').Value = '1.00'
$ws.Range('E46').Value = '  +0.04%  '

$ws.Range('D47').Value = '66.05'
$ws.Range('E47').Value = '  +3.39%  '

$ws.Range('E48').Value = '  +1.47%  '

$ws.Range('D49').Value = '1.750.37'
$ws.Range('E49').Value = '  +3.13%  '

$ws.Range('D50').Value = '0.875'
$ws.Range('E50').Value = '  +1.15%  '

$ws.Range('D51').Value = '86.75'
$ws.Range('E51').Value = '  +2.19%  '

# Restore default (unstyled) cell style on column D now that the text values
# are locked in, so only the cell contents differ from the original workbook.
$ws.Range('D2:D51').Style = 'Normal'
